# Fuel Prod Imp Exp Balancing Priorities.xlsx - apply author's edits
# (updated 4.0 files and mdl)

$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

# --- About sheet: bump the "last updated" date stamp (C1) ---------------
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: re-prioritize "hard coal" production/imports/exports -
# production(B3): 3 -> 1 ; imports(C3): 2 -> 3 ; exports(D3): 1 -> 2
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# --- FPIEBP sheet: move the active selection from F4 to E3 --------------
$wsFPIEBP.Activate()
[void]$wsFPIEBP.Range("E3").Select()
